$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Magnesium chloride unit price - update lower/upper bounds to plain
# values (no longer derived from the midpoint via formula).
$ws.Range("E8").Value = 0.38
$ws.Range("G8").Value = 0.349
$ws.Range("I8").Value = 0.411

# Row 9: Zinc sulfate unit price - same treatment.
$ws.Range("E9").Value = 0.795
$ws.Range("G9").Value = 0.657
$ws.Range("I9").Value = 0.931

# Reflect the new active selection left by the editor (rows 8:9 selected).
[void]$ws.Range("A8:XFD9").Select()
